$d = $word.ActiveDocument
$styles = $d.Styles

# docDefaults rPrDefault/rPr/rFonts@eastAsia and the Normal style's
# rPr/rFonts@eastAsia both resolve through Word's "Normal" style font
# (Normal carries the document default run properties); update eastAsia
# (w:eastAsia) from "DejaVu Sans" to "Tahoma".
$normal = $styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

# Heading style: same eastAsia swap.
$heading = $styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# List style previously had an empty rPr; it now carries an explicit
# complex-script font override (w:cs) of "DejaVu Sans".
$list = $styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

# Caption style gains the same explicit complex-script font override.
$caption = $styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

# Index style previously had an empty rPr; same complex-script override.
$index = $styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
